# TradingModel_v2 - 2021/11/19 data updated
# Append the 2021-11-19 trading rows (13-19) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 13; Date = 44519; StockId = 3122; ProfitPercent = -0.02;               ProfitMoney = -141 },
    @{ Row = 14; Date = 44519; StockId = 3221; ProfitPercent = -0.01;               ProfitMoney = -21.34999999999991 },
    @{ Row = 15; Date = 44519; StockId = 6282; ProfitPercent = -0.07000000000000001; ProfitMoney = -387.5 },
    @{ Row = 16; Date = 44519; StockId = 6196; ProfitPercent = -0.04;               ProfitMoney = -245 },
    @{ Row = 17; Date = 44519; StockId = 3035; ProfitPercent = 0.07000000000000001;  ProfitMoney = 432 },
    @{ Row = 18; Date = 44519; StockId = 6411; ProfitPercent = 0.18;                ProfitMoney = 1105 },
    @{ Row = 19; Date = 44519; StockId = 3141; ProfitPercent = 0.19;                ProfitMoney = 1093.5 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Date column keeps the same date/time number format + style as the rows above it.
    $ws.Range("A$row").Value = $r.Date
    $ws.Range("A$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B$row").Value = $r.StockId
    $ws.Range("C$row").Value = $r.ProfitPercent
    $ws.Range("D$row").Value = $r.ProfitMoney

    # Columns E:K stay blank (inline empty strings), matching the rest of the table.
    $ws.Range("E12:K12").Copy($ws.Range("E$row"))
}
